$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.8080601414048771
$ws.Range("E2").Value = 0.8080601414048771

# Row 3
$ws.Range("D3").Value = 0.001066615700799473
$ws.Range("E3").Value = 0.001066615700799473

# Row 4
$ws.Range("D4").Value = 0.9841227757605735
$ws.Range("E4").Value = 0.9841227757605735

# Row 5
$ws.Range("D5").Value = 0.001973751605342619
$ws.Range("E5").Value = 0.001973751605342619

# Row 6
$ws.Range("D6").Value = 0.2928436514589686
$ws.Range("E6").Value = 0.2928436514589686

# Row 7
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0

# Row 8
$ws.Range("D8").Value = 0.9615572033427451
$ws.Range("E8").Value = 0.03844279665725492

# Row 9
$ws.Range("D9").Value = 0.9999887454760186
$ws.Range("E9").Value = 0.0000112545239814299

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.3108077500029213
$ws.Range("E10").Value = 0.6891922499970786

# Row 11
$ws.Range("D11").Value = 0.9246098342588457
$ws.Range("E11").Value = 0.07539016574115431
$ws.Range("F11").Value = 0.7429165244102478
$ws.Range("G11").Value = 0.7

# Row 12
$ws.Range("D12").Value = 0.9631034468299031
$ws.Range("E12").Value = 0.9631034468299031

# Row 13
$ws.Range("D13").Value = 0.00003452793200000132
$ws.Range("E13").Value = 0.00003452793200000132

# Row 14
$ws.Range("D14").Value = 0.9989885255933324
$ws.Range("E14").Value = 0.9989885255933324

# Row 15
$ws.Range("D15").Value = 0.00004432059862822761
$ws.Range("E15").Value = 0.00004432059862822761

# Row 16
$ws.Range("D16").Value = 0.112315230563792
$ws.Range("E16").Value = 0.112315230563792

# Row 17
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0

# Row 18
$ws.Range("D18").Value = 0.8680781779308601
$ws.Range("E18").Value = 0.1319218220691399

# Row 19
$ws.Range("D19").Value = 0.999999999596801
$ws.Range("E19").Value = 0.0000000004031990297193033

# Row 20
$ws.Range("D20").Value = 0.7462945916079222
$ws.Range("E20").Value = 0.2537054083920778

# Row 21
$ws.Range("D21").Value = 0.9208224608853164
$ws.Range("E21").Value = 0.07917753911468361
$ws.Range("F21").Value = 1.083179831504822
$ws.Range("G21").Value = 0.8

$wb.Save()
